$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume (E) columns for rows 2-51 so that
# numeric-looking strings (e.g. "0.9968", "30.82") are preserved exactly as text
# rather than being auto-converted to numbers (which would drop trailing zeros).
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

# --- Apply cell updates from the source diff ---
$ws.Range("D2").Value = "26.203.17"
$ws.Range("E2").Value = "  +4.65%  "
$ws.Range("D3").Value = "1.693.83"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("D4").Value = "0.9968"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "240.08"
$ws.Range("E5").Value = "  +3.28%  "
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "0.4686"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").Value = "0.2639"
$ws.Range("D9").Value = "0.06194"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").Value = "1.682.19"
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("D11").Value = "0.07059"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "15.23"
$ws.Range("E12").Value = "  +6.55%  "
$ws.Range("D13").Value = "4.416"
$ws.Range("E13").Value = "  +2.78%  "
$ws.Range("D14").Value = "0.5879"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("D15").Value = "75.97"
$ws.Range("E15").Value = "  +3.64%  "
$ws.Range("D16").Value = "0.9988"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "0.9974"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "26.162.98"
$ws.Range("E18").Value = "  +4.56%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "11.59"
$ws.Range("E19").Value = "  +3.08%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.000006779"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("D21").Value = "1.898.62"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("D22").Value = "4.539"
$ws.Range("E22").Value = "  +6.21%  "
$ws.Range("D23").Value = "8.793"
$ws.Range("E23").Value = "  +4.23%  "
$ws.Range("D24").Value = "5.301"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "134.49"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").Value = "15.11"
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("D27").Value = "1.396"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").Value = "1.742"
$ws.Range("E28").Value = "  +7.03%  "
$ws.Range("D29").Value = "105.64"
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("D30").Value = "3.995"
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("D31").Value = "3.680"
$ws.Range("E31").Value = "  +5.00%  "
$ws.Range("D32").Value = "0.07768"
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("D33").Value = "0.04388"
$ws.Range("E33").Value = "  +3.80%  "
$ws.Range("D34").Value = "2.618"
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "0.9687"
$ws.Range("E35").Value = "  +4.26%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6168"
$ws.Range("E36").Value = "  +5.34%  "
$ws.Range("D37").Value = "0.9238"
$ws.Range("E37").Value = "  +6.98%  "
$ws.Range("D38").Value = "110.70"
$ws.Range("E38").Value = "  +13.16%  "
$ws.Range("D39").Value = "2.382"
$ws.Range("E39").Value = "  -7.20%  "
$ws.Range("D40").Value = "0.9998"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").Value = "1.891"
$ws.Range("E41").Value = "  +6.96%  "
$ws.Range("D42").Value = "0.01468"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").Value = "0.3777"
$ws.Range("E43").Value = "  +3.34%  "
$ws.Range("D44").Value = "5.108"
$ws.Range("E44").Value = "  +10.53%  "
$ws.Range("D45").Value = "0.1133"
$ws.Range("E45").Value = "  +3.99%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "6.227"
$ws.Range("E46").Value = "  +3.14%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05320"
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("D48").Value = "30.82"
$ws.Range("E48").Value = "  +7.62%  "
$ws.Range("D49").Value = "7.671"
$ws.Range("E49").Value = "  +8.05%  "
$ws.Range("D50").Value = "1.216"
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.3354"
$ws.Range("E51").Value = "  +3.01%  "

# Restore default (Normal) style on the Price/Volume range so no stray number-format
# styling remains attached to these cells (matches original unstyled cells).
$priceVolRange.Style = "Normal"

